# Added v0 scripts for western blot 2v0, 3v0 to test change lysis pull to
# waste syringe first.
#
# New rows are appended to the "v0" worksheet describing two new scripts
# (v0_script_2v0 / "Western Blot 1" and v0_script_3v0 / "Western Blot 1 -
# change lysis pull to waste syringe"), with the RIPA incubation-time note
# highlighted in yellow, plus a long explanatory note on row 6.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("v0")

# Extend the header row with an extra (empty, bold) cell in column N -
# matches the used-range growing out to column N.
$ws2.Range("N3").Value = ""
$ws2.Range("N3").Font.Bold = $true

# Seed brand-new shared strings in the same order they are first
# encountered left-to-right, top-to-bottom across the two new rows, so the
# shared string table lines up exactly.
$ws2.Range("A5").Value = "v0_script_2v0"
$ws2.Range("B5").Value = "Western Blot 1"
$ws2.Range("B6").Value = "Western Blot 1 - change lysis pull to waste syringe"
$ws2.Range("A6").Value = "v0_script_3v0"
$ws2.Range("J5").Value = "5 mins (RIPA - 500 uL)"
$ws2.Range("L6").Value = "***CHANGE FROM 1V0: after add RIPA lysis buffer, pull to WASTE syringe first, 5 min incubation, then LAST pull to LYSATE syringe (test with dyes, may need to increase pull to 700 uL to match QIAzol since this gets lysis buffer to outlet of chip, but may leave reservoir empty, affect last pull?)"

# Remaining row 5 cells (all reuse existing shared strings).
$ws2.Range("C5").Value = "5 mL"
$ws2.Range("D5").Value = "5 mL"
$ws2.Range("E5").Value = "1 hour"
$ws2.Range("F5").Value = "0.5 mL"
$ws2.Range("G5").Value = "15 mL/hr"
$ws2.Range("H5").Value = "15 mL/hr"
$ws2.Range("I5").Value = "200-800-1000"
$ws2.Range("K5").Value = "N"

# Remaining row 6 cells (all reuse existing shared strings).
$ws2.Range("C6").Value = "5 mL"
$ws2.Range("D6").Value = "5 mL"
$ws2.Range("E6").Value = "1 hour"
$ws2.Range("F6").Value = "0.5 mL"
$ws2.Range("G6").Value = "15 mL/hr"
$ws2.Range("H6").Value = "15 mL/hr"
$ws2.Range("I6").Value = "200-800-1000"
$ws2.Range("J6").Value = "5 mins (RIPA - 500 uL)"
$ws2.Range("K6").Value = "N"

# Highlight the RIPA incubation note (J5, J6) and the long change note (L6)
# in yellow.
$ws2.Range("J5").Interior.Color = 65535
$ws2.Range("J6").Interior.Color = 65535
$ws2.Range("L6").Interior.Color = 65535

# Switch the active tab from "r0" to "v0", scroll so column F is rendered
# at the sheet's left edge, and leave the selection on I9.
$ws2.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws2.Range("I9").Select()
